$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1611941447041101
$ws.Range("E2").Value = 11.07446084772978
$ws.Range("F2").Value = 32.11428885690793
